# Update phone numbers in column C (rows 2-5): 87824121996 -> 87824121997
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 87824121997
$ws.Range("C3").Value = 87824121997
$ws.Range("C4").Value = 87824121997
$ws.Range("C5").Value = 87824121997

# Update the active selection from B8 to D6
$ws.Range("D6").Select()
